$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

$newRow = 66

# Apply the highlighted-row formatting (fill + border + number formats) used
# for marked rows onto the newly appended row *before* writing values, so the
# date cell keeps its numeric date format instead of falling back to text.
$ws.Range("A42:J42").Copy()
$ws.Range("A66:J66").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Append new data row for 15. 5. 2020 (CO municipality-level data)
$ws.Cells.Item($newRow, 1).Value = 43966
$ws.Cells.Item($newRow, 2).Value = 68852
$ws.Cells.Item($newRow, 3).Value = 1151
$ws.Cells.Item($newRow, 4).Value = 1465
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 27
$ws.Cells.Item($newRow, 7).Value = 6
$ws.Cells.Item($newRow, 8).Value = 2
$ws.Cells.Item($newRow, 9).Value = 103
$ws.Cells.Item($newRow, 10).Value = 0

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Tabela1")
$table.Resize($ws.Range("A1:J66"))

# Update selection / active cell to mirror the workbook view state
$ws.Range("A66:J66").Select()
